$d = $word.ActiveDocument

# 1. Add trailing "፣" after "ፊል" in Micky's question
$d.Content.Find.Execute(
    "ሚኪ እጅግ ተገርሞ፣ ጠየቀው፡ ስንት አመታቸው ነወ? ፊል", $true, $false, $false, $false, $false,
    $true, 1, $false, "ሚኪ እጅግ ተገርሞ፣ ጠየቀው፡ ስንት አመታቸው ነወ? ፊል፣", 2)

# 2. Add "-" after "፡" in "እሰጥሃለሁ፡ አነተ"
$d.Content.Find.Execute(
    "“ትነግርኛለህ! ፍንጭ እሰጥሃለሁ፡ አነተ", $true, $false, $false, $false, $false,
    $true, 1, $false, "“ትነግርኛለህ! ፍንጭ እሰጥሃለሁ፡- አነተ", 2)

# 3. Remove "ልጆች " from "የሶስቱን ልጆች እድሜ"
$d.Content.Find.Execute(
    "የሶስቱን ልጆች እድሜ አንድ ላይ ስታባዛው", $true, $false, $false, $false, $false,
    $true, 1, $false, "የሶስቱን እድሜ አንድ ላይ ስታባዛው", 2)

# 4. Change "፡" to "፡--  " (two hyphens, two spaces) in "እንድህ አለው፡ አዝናለሁ"
$d.Content.Find.Execute(
    "እንድህ አለው፡ አዝናለሁ ፊል፣ ነገር ግን", $true, $false, $false, $false, $false,
    $true, 1, $false, "እንድህ አለው፡--  አዝናለሁ ፊል፣ ነገር ግን", 2)

# 5. Add "ቢያንስ " at the beginning of "አንድ ተጨማሪ እፈልጋለሁ..."
$d.Content.Find.Execute(
    "አንድ ተጨማሪ እፈልጋለሁ፡፡” ፊል እንድህ አለ፡- “አወ", $true, $false, $false, $false, $false,
    $true, 1, $false, "ቢያንስ አንድ ተጨማሪ እፈልጋለሁ፡፡” ፊል እንድህ አለ፡- “አወ", 2)
